$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PBIReports")

$ws.Range("H2:H7").Value = "Finance App"

$ws.Range("I2").Value = "Financial Overview"
$ws.Range("I3").Value = "Income Statement by Month"
$ws.Range("I4").Value = "Balance Sheet by Month"
$ws.Range("I5").Value = "Budget Comparison"
$ws.Range("I6").Value = "Liquidity KPIs"
$ws.Range("I7").Value = "Profitability"
